$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header suffix to "_FV2404" and the "_new" header suffix to "_FV2410"
$oldNames = @(
    "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
    "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old"
)
$newNamesFV2404 = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newNamesFV2404[$i]
}

$oldNamesNew = @(
    "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
    "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)
$newNamesFV2410 = @(
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $oldNamesNew.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newNamesFV2410[$i]
}

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into a formatted Excel Table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null
